# Commit: "some style changes, and syllabuses for summer courses"
#
# For this document (word/research-seminar.docx) the portion of the
# change that is reachable through the Word object model is the
# addition of the built-in-looking "Footnote Text" paragraph style to
# the style sheet (word/styles.xml). It is defined exactly like Word's
# normal "Footnote Text" quick style: based on Normal, followed by
# itself, uiPriority 9, unhidden-when-used, and part of the Quick
# Style gallery.
#
# (The nsid GUIDs inside word/numbering.xml that also changed in the
# source diff are an internal bookkeeping identifier for each list
# definition that Word has never exposed anywhere in its object model
# -- there is no property on ListTemplate/ListFormat/List, no Find/
# Replace target, nothing -- so it is not something that can be driven
# from COM automation; it is left untouched here.)

$d = $word.ActiveDocument

$footnoteText = $d.Styles.Add("Footnote Text", 1)
$footnoteText.BaseStyle = "Normal"
$footnoteText.NextParagraphStyle = "FootnoteText"
$footnoteText.Priority = 9
$footnoteText.UnhideWhenUsed = $true
$footnoteText.QuickStyle = $true

Write-Output "Added FootnoteText style"
